# The sheet holds a daily price log for "Cebollín" at Feria Lagunitas de
# Puerto Montt, sorted with the most recent reading first. A new weekly
# reading was added at the top of the data block (row 226, just below the
# header row), pushing every existing record down by one row (226->227,
# 227->228, ... 316->317).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing data block (rows 226:316) down by one row. Excel's
# native row-insert shifts the cell contents (and copies formatting, e.g.
# the date number format on column D) down automatically.
$ws.Rows("226:226").Insert()

# Populate the newly opened row 226 with the new reading.
$ws.Range("A226").Value = 4
$ws.Range("B226").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C226").Value = "Los Lagos"
$ws.Range("D226").Value = 44784
$ws.Range("E226").Value = 10
$ws.Range("F226").Value = 100112037
$ws.Range("G226").Value = "Cebollín"
$ws.Range("H226").Value = "Sin especificar"
$ws.Range("I226").Value = "Segunda"
$ws.Range("J226").Value = 70
$ws.Range("K226").Value = 9500
$ws.Range("L226").Value = 9500
$ws.Range("M226").Value = 9500
$ws.Range("N226").Value = "$/paquete 36 unidades"
$ws.Range("O226").Value = "Región Metropolitana"
$ws.Range("P226").Value = 264
$ws.Range("Q226").Value = 36
$ws.Range("R226").Value = "Hortaliza"
